$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("219×5=", $true, $false, $false, $false, $false, $true, 1, $false, "633×4=", 2) | Out-Null
$d.Content.Find.Execute("214×8=", $true, $false, $false, $false, $false, $true, 1, $false, "185×7=", 2) | Out-Null
$d.Content.Find.Execute("636×3=", $true, $false, $false, $false, $false, $true, 1, $false, "818×7=", 2) | Out-Null
$d.Content.Find.Execute("464×3=", $true, $false, $false, $false, $false, $true, 1, $false, "959×6=", 2) | Out-Null
$d.Content.Find.Execute("375×6=", $true, $false, $false, $false, $false, $true, 1, $false, "913×6=", 2) | Out-Null
$d.Content.Find.Execute("219×2=", $true, $false, $false, $false, $false, $true, 1, $false, "482×8=", 2) | Out-Null
$d.Content.Find.Execute("329×4=", $true, $false, $false, $false, $false, $true, 1, $false, "839×2=", 2) | Out-Null
$d.Content.Find.Execute("865×2=", $true, $false, $false, $false, $false, $true, 1, $false, "263×5=", 2) | Out-Null
$d.Content.Find.Execute("767×2=", $true, $false, $false, $false, $false, $true, 1, $false, "716×6=", 2) | Out-Null
$d.Content.Find.Execute("811×2=", $true, $false, $false, $false, $false, $true, 1, $false, "596×5=", 2) | Out-Null
$d.Content.Find.Execute("293×2=", $true, $false, $false, $false, $false, $true, 1, $false, "784×4=", 2) | Out-Null
$d.Content.Find.Execute("846×5=", $true, $false, $false, $false, $false, $true, 1, $false, "409×7=", 2) | Out-Null
$d.Content.Find.Execute("926×4=", $true, $false, $false, $false, $false, $true, 1, $false, "668×8=", 2) | Out-Null
$d.Content.Find.Execute("362×6=", $true, $false, $false, $false, $false, $true, 1, $false, "149×3=", 2) | Out-Null
$d.Content.Find.Execute("600×3=", $true, $false, $false, $false, $false, $true, 1, $false, "484×7=", 2) | Out-Null
$d.Content.Find.Execute("975×8=", $true, $false, $false, $false, $false, $true, 1, $false, "962×7=", 2) | Out-Null
$d.Content.Find.Execute("212×9=", $true, $false, $false, $false, $false, $true, 1, $false, "647×9=", 2) | Out-Null
$d.Content.Find.Execute("960×8=", $true, $false, $false, $false, $false, $true, 1, $false, "739×3=", 2) | Out-Null
$d.Content.Find.Execute("279×3=", $true, $false, $false, $false, $false, $true, 1, $false, "858×5=", 2) | Out-Null
$d.Content.Find.Execute("786×8=", $true, $false, $false, $false, $false, $true, 1, $false, "834×5=", 2) | Out-Null
$d.Content.Find.Execute("751×7=", $true, $false, $false, $false, $false, $true, 1, $false, "854×4=", 2) | Out-Null
$d.Content.Find.Execute("381×7=", $true, $false, $false, $false, $false, $true, 1, $false, "761×3=", 2) | Out-Null
$d.Content.Find.Execute("670×9=", $true, $false, $false, $false, $false, $true, 1, $false, "899×8=", 2) | Out-Null
$d.Content.Find.Execute("582×9=", $true, $false, $false, $false, $false, $true, 1, $false, "785×4=", 2) | Out-Null
$d.Content.Find.Execute("232×4=", $true, $false, $false, $false, $false, $true, 1, $false, "258×2=", 2) | Out-Null
